$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F header: same style as the other header cells (A1:E1) ---
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Convert column A (rows 2-7) from text timestamps to real date-time numbers ---
$ws.Range("A2").Value = 45685.64979224537
$ws.Range("A3").Value = 45685.64983969907
$ws.Range("A4").Value = 45685.66042303241
$ws.Range("A5").Value = 45685.64978993055
$ws.Range("A6").Value = 45685.64983738426
$ws.Range("A7").Value = 45685.66042071759

# --- New rows 8-13 ---
$ws.Range("A8").Value = 45685.67012094907
$ws.Range("B8").Value = 3013.4
$ws.Range("C8").Value = 11.91
$ws.Range("D8").Value = 3.532313176563807
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Mała Gra"

$ws.Range("A9").Value = 45685.67027604167
$ws.Range("B9").Value = 3026.8
$ws.Range("C9").Value = 11.71
$ws.Range("D9").Value = 3.489612885883876
$ws.Range("E9").Value = "10-15"
$ws.Range("F9").Value = "Mała Gra"

$ws.Range("A10").Value = 45685.67605960649
$ws.Range("B10").Value = 3526.5
$ws.Range("C10").Value = 14.84
$ws.Range("D10").Value = 3.842985357557024
$ws.Range("E10").Value = "10-15"
$ws.Range("F10").Value = "Mała Gra"

$ws.Range("A11").Value = 45685.66946122685
$ws.Range("B11").Value = 2956.4
$ws.Range("C11").Value = 8.94
$ws.Range("D11").Value = 3.075016839163645
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "Mała Gra"

$ws.Range("A12").Value = 45685.67027372685
$ws.Range("B12").Value = 3026.6
$ws.Range("C12").Value = 8.42
$ws.Range("D12").Value = 2.981690219470433
$ws.Range("E12").Value = "5-10"
$ws.Range("F12").Value = "Mała Gra"

$ws.Range("A13").Value = 45685.68411863426
$ws.Range("B13").Value = 4222.8
$ws.Range("C13").Value = 9.789999999999999
$ws.Range("D13").Value = 2.995159932545254
$ws.Range("E13").Value = "5-10"
$ws.Range("F13").Value = "Mała Gra"

# --- Column F for existing rows (2-7): "Duża Gra" ---
$ws.Range("F2").Value = "Duża Gra"
$ws.Range("F3").Value = "Duża Gra"
$ws.Range("F4").Value = "Duża Gra"
$ws.Range("F5").Value = "Duża Gra"
$ws.Range("F6").Value = "Duża Gra"
$ws.Range("F7").Value = "Duża Gra"

# --- Apply the date/time number format. Establish it once (creating the
#     164 -> 165 numFmt pair, matching Excel's own intermediate-registration
#     behaviour), then propagate the final format to every timestamp cell. ---
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$dateFmt = $ws.Range("A2").NumberFormat
$ws.Range("A3").NumberFormat = $dateFmt
$ws.Range("A4").NumberFormat = $dateFmt
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("A6").NumberFormat = $dateFmt
$ws.Range("A7").NumberFormat = $dateFmt
$ws.Range("A8").NumberFormat = $dateFmt
$ws.Range("A9").NumberFormat = $dateFmt
$ws.Range("A10").NumberFormat = $dateFmt
$ws.Range("A11").NumberFormat = $dateFmt
$ws.Range("A12").NumberFormat = $dateFmt
$ws.Range("A13").NumberFormat = $dateFmt
